$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as TEXT (matches the source data, which is all
# scraped/formatted text such as "0.2600" or "24.989.99", not numerics).
# Briefly switching the cell to a text format forces Excel to keep the
# literal string instead of auto-coercing it to a Number/Date, then we
# restore "General" so the cell keeps its original (unformatted) look.
function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.NumberFormat = "General"
}

Set-TextValue "D2" "24.989.99"
$ws.Range("E2").Value = "  -3.87%  "

Set-TextValue "D3" "1.642.37"
$ws.Range("E3").Value = "  -5.76%  "

Set-TextValue "D4" "0.9993"
$ws.Range("E4").Value = "  -0.06%  "

Set-TextValue "D5" "233.07"
$ws.Range("E5").Value = "  -5.86%  "

Set-TextValue "D7" "0.4759"
$ws.Range("E7").Value = "  -5.58%  "

Set-TextValue "D8" "0.2600"

Set-TextValue "D9" "0.06113"
$ws.Range("E9").Value = "  -1.27%  "

Set-TextValue "D10" "0.07023"
$ws.Range("E10").Value = "  -3.11%  "

Set-TextValue "D11" "1.648.31"
$ws.Range("E11").Value = "  -5.41%  "

Set-TextValue "D12" "14.57"
$ws.Range("E12").Value = "  -3.69%  "

Set-TextValue "D13" "0.5891"
$ws.Range("E13").Value = "  -10.10%  "

Set-TextValue "D14" "4.337"
$ws.Range("E14").Value = "  -7.50%  "

Set-TextValue "D15" "73.63"
$ws.Range("E15").Value = "  -5.14%  "

Set-TextValue "D16" "1.000"
$ws.Range("E16").Value = "  +0.04%  "

$ws.Range("E17").Value = "  -0.01%  "

Set-TextValue "D18" "24.994.63"
$ws.Range("E18").Value = "  -3.93%  "

Set-TextValue "D19" "0.000006595"
$ws.Range("E19").Value = "  -4.06%  "

Set-TextValue "D20" "11.22"
$ws.Range("E20").Value = "  -6.10%  "

Set-TextValue "D21" "1.857.93"
$ws.Range("E21").Value = "  -5.61%  "

Set-TextValue "D22" "4.284"
$ws.Range("E22").Value = "  -4.63%  "

Set-TextValue "D23" "8.555"
$ws.Range("E23").Value = "  -1.95%  "

Set-TextValue "D24" "5.238"
$ws.Range("E24").Value = "  -3.17%  "

Set-TextValue "D25" "133.75"
$ws.Range("E25").Value = "  -1.53%  "

Set-TextValue "D26" "14.91"
$ws.Range("E26").Value = "  -2.32%  "

$ws.Range("E27").Value = "  -7.50%  "

Set-TextValue "D28" "103.40"
$ws.Range("E28").Value = "  -1.99%  "

Set-TextValue "D29" "1.634"
$ws.Range("E29").Value = "  -8.75%  "

Set-TextValue "D30" "3.892"
$ws.Range("E30").Value = "  -1.80%  "

Set-TextValue "D31" "0.07658"
$ws.Range("E31").Value = "  -6.29%  "

Set-TextValue "D32" "3.583"
$ws.Range("E32").Value = "  -2.66%  "

$ws.Range("E33").Value = "  +0.08%  "

Set-TextValue "D34" "0.04285"
$ws.Range("E34").Value = "  -8.74%  "

Set-TextValue "D35" "2.573"
$ws.Range("E35").Value = "  -3.26%  "

$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D36" "0.9263"
$ws.Range("E36").Value = "  -6.93%  "

$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D37" "0.5895"
$ws.Range("E37").Value = "  -3.73%  "

Set-TextValue "D38" "2.581"
$ws.Range("E38").Value = "  -6.26%  "

Set-TextValue "D39" "0.8747"
$ws.Range("E39").Value = "  +7.85%  "

Set-TextValue "D40" "0.9999"
$ws.Range("E40").Value = "  +0.00%  "

Set-TextValue "D41" "0.01509"
$ws.Range("E41").Value = "  -7.11%  "

Set-TextValue "D42" "98.71"
$ws.Range("E42").Value = "  -2.31%  "

Set-TextValue "D43" "1.760"
$ws.Range("E43").Value = "  -8.97%  "

Set-TextValue "D44" "0.3693"
$ws.Range("E44").Value = "  -5.86%  "

Set-TextValue "D45" "4.669"
$ws.Range("E45").Value = "  -6.91%  "

Set-TextValue "D46" "0.1099"
$ws.Range("E46").Value = "  -5.38%  "

$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D47" "0.05211"
$ws.Range("E47").Value = "  -1.59%  "

$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D48" "6.083"
$ws.Range("E48").Value = "  -3.96%  "

Set-TextValue "D49" "28.92"
$ws.Range("E49").Value = "  -6.02%  "

Set-TextValue "D51" "0.9995"
$ws.Range("E51").Value = "  +0.19%  "
